$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly roll-forward of the price table: each date block's values shift to
# the next reporting week, with the Fecha / Calidad / Volumen / Precio
# mínimo / Precio máximo / Precio promedio ponderado / Precio $/Kg columns
# updated accordingly.

$ws.Range("D5").Value  = 44322
$ws.Range("M5").Value  = 200
$ws.Range("N5").Value  = 7000
$ws.Range("O5").Value  = 7500
$ws.Range("P5").Value  = 7250
$ws.Range("S5").Value  = 2417

$ws.Range("D6").Value  = 44322
$ws.Range("N6").Value  = 6000
$ws.Range("O6").Value  = 6500
$ws.Range("P6").Value  = 6250
$ws.Range("S6").Value  = 2083

$ws.Range("D7").Value  = 44322
$ws.Range("M7").Value  = 100

$ws.Range("L8").Value  = "Especial"
$ws.Range("M8").Value  = 100
$ws.Range("N8").Value  = 6500
$ws.Range("O8").Value  = 7000
$ws.Range("P8").Value  = 6750
$ws.Range("S8").Value  = 2250

$ws.Range("D9").Value  = 44172
$ws.Range("L9").Value  = "Primera"
$ws.Range("N9").Value  = 5500
$ws.Range("O9").Value  = 6000
$ws.Range("P9").Value  = 5750
$ws.Range("S9").Value  = 1917

$ws.Range("D10").Value = 44172
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 160
$ws.Range("N10").Value = 5000
$ws.Range("O10").Value = 5500
$ws.Range("P10").Value = 5250
$ws.Range("S10").Value = 1750

$ws.Range("D11").Value = 44172
$ws.Range("L11").Value = "Tercera"
$ws.Range("M11").Value = 140
$ws.Range("N11").Value = 3500
$ws.Range("O11").Value = 4000
$ws.Range("P11").Value = 3750
$ws.Range("S11").Value = 1250

$ws.Range("D15").Value = 44249
$ws.Range("N15").Value = 6000
$ws.Range("O15").Value = 7000
$ws.Range("P15").Value = 6500
$ws.Range("S15").Value = 2167

$ws.Range("D16").Value = 44249
$ws.Range("N16").Value = 4500
$ws.Range("O16").Value = 5000
$ws.Range("P16").Value = 4750
$ws.Range("S16").Value = 1583

$ws.Range("D17").Value = 44351
$ws.Range("L17").Value = "Especial"
$ws.Range("M17").Value = 160
$ws.Range("N17").Value = 7500
$ws.Range("O17").Value = 8000
$ws.Range("P17").Value = 7750
$ws.Range("S17").Value = 2583

$ws.Range("D18").Value = 44351
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 6000
$ws.Range("O18").Value = 6500
$ws.Range("P18").Value = 6250
$ws.Range("S18").Value = 2083

$ws.Range("D19").Value = 44351
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 200
$ws.Range("N19").Value = 4500
$ws.Range("O19").Value = 5000
$ws.Range("P19").Value = 4750
$ws.Range("S19").Value = 1583

$ws.Range("D20").Value = 44334
$ws.Range("L20").Value = "Especial"
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 7000
$ws.Range("O20").Value = 8000
$ws.Range("P20").Value = 7500
$ws.Range("S20").Value = 2500

$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 160
$ws.Range("N21").Value = 6000
$ws.Range("O21").Value = 7000
$ws.Range("P21").Value = 6500
$ws.Range("S21").Value = 2167

$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 120

$ws.Range("L23").Value = "Tercera"
$ws.Range("M23").Value = 70
$ws.Range("N23").Value = 3500
$ws.Range("O23").Value = 4000
$ws.Range("P23").Value = 3750
$ws.Range("S23").Value = 1250

$ws.Range("D24").Value = 44200
$ws.Range("L24").Value = "Especial"
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 4500
$ws.Range("O24").Value = 5000
$ws.Range("P24").Value = 4750
$ws.Range("S24").Value = 1583

$ws.Range("D25").Value = 44200
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 80
$ws.Range("N25").Value = 3500
$ws.Range("O25").Value = 4000
$ws.Range("P25").Value = 3750
$ws.Range("S25").Value = 1250

$ws.Range("D26").Value = 44200
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 120
$ws.Range("N26").Value = 2500
$ws.Range("O26").Value = 3000
$ws.Range("P26").Value = 2750
$ws.Range("S26").Value = 917
